$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion message text (A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.49 = 17660.45 pesos`n✅ 17660.45 pesos = 4.44 = 896.33 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 222.872
$ws2.Range("O10").Value = 3936.02
$ws2.Range("N12").Value = 3982
$ws2.Range("O12").Value = 202.1
